$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Notes": update the specific-issue description text
# ---------------------------------------------------------------
$wsNotes = $wb.Worksheets.Item("Notes")
$wsNotes.Range("A3").Value = "Specific issue: total_num is not identical after grouping by survey-gene-locus"

# ---------------------------------------------------------------
# Sheet "studies": lower-case the study_ID header -> study_id
# ---------------------------------------------------------------
$wsStudies = $wb.Worksheets.Item("studies")
$wsStudies.Range("A1").Value = "study_id"
$wsStudies.Activate()
$wsStudies.Range("A2").Select()

# ---------------------------------------------------------------
# Sheet "surveys": lower-case survey_ID -> survey_id, rename
# lat/lon -> latitude/longitude, and apply new header styling
# ---------------------------------------------------------------
$wsSurveys = $wb.Worksheets.Item("surveys")
$wsSurveys.Range("B1").Value = "survey_id"
$wsSurveys.Range("E1").Value = "latitude"
$wsSurveys.Range("F1").Value = "longitude"

$wsSurveys.Range("A1:G1").Font.Color = 0
$wsSurveys.Range("K1").Font.Color = 0
$wsSurveys.Range("H1:J1").Font.Color = 0
$wsSurveys.Range("H1:J1").NumberFormat = "@"

$wsSurveys.Activate()
$wsSurveys.Range("A1:K1").Select()

# ---------------------------------------------------------------
# Sheet "counts": update variant strings and numbers
# ---------------------------------------------------------------
$wsCounts = $wb.Worksheets.Item("counts")
$wsCounts.Range("B2").Value = "crt:1:A"
$wsCounts.Range("C2").Value = 2
$wsCounts.Range("B3").Value = "crt:1:C"
$wsCounts.Range("C3").Value = 2
$wsCounts.Range("D3").Value = 11

$wsCounts.Activate()
$wsCounts.Range("D5").Select()
